{"js": "// Replace the 25 two-digit-divided-by-one-digit practice answers in the\n// table with the newly generated values, matching each unique original\n// string and swapping it for its replacement.\nconst replacements = [\n  [\"64\u00f75=12, 4\", \"67\u00f79=7, 4\"],\n  [\"57\u00f72=28, 1\", \"45\u00f72=22, 1\"],\n  [\"65\u00f79=7, 2\", \"88\u00f79=9, 7\"],\n  [\"62\u00f79=6, 8\", \"50\u00f78=6, 2\"],\n  [\"92\u00f74=23, 0\", \"23\u00f76=3, 5\"],\n  [\"60\u00f74=15, 0\", \"51\u00f77=7, 2\"],\n  [\"78\u00f77=11, 1\", \"83\u00f76=13, 5\"],\n  [\"79\u00f77=11, 2\", \"38\u00f72=19, 0\"],\n  [\"49\u00f73=16, 1\", \"19\u00f72=9, 1\"],\n  [\"78\u00f79=8, 6\", \"28\u00f73=9, 1\"],\n  [\"93\u00f78=11, 5\", \"81\u00f77=11, 4\"],\n  [\"76\u00f75=15, 1\", \"49\u00f77=7, 0\"],\n  [\"45\u00f75=9, 0\", \"29\u00f75=5, 4\"],\n  [\"36\u00f78=4, 4\", \"97\u00f79=10, 7\"],\n  [\"25\u00f78=3, 1\", \"76\u00f72=38, 0\"],\n  [\"66\u00f74=16, 2\", \"46\u00f76=7, 4\"],\n  [\"42\u00f74=10, 2\", \"32\u00f73=10, 2\"],\n  [\"11\u00f75=2, 1\", \"45\u00f78=5, 5\"],\n  [\"67\u00f75=13, 2\", \"55\u00f78=6, 7\"],\n  [\"48\u00f78=6, 0\", \"78\u00f78=9, 6\"],\n  [\"43\u00f76=7, 1\", \"40\u00f79=4, 4\"],\n  [\"68\u00f72=34, 0\", \"66\u00f75=13, 1\"],\n  [\"69\u00f72=34, 1\", \"52\u00f77=7, 3\"],\n  [\"62\u00f72=31, 0\", \"61\u00f79=6, 7\"],\n  [\"29\u00f78=3, 5\", \"18\u00f75=3, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 two-digit-divided-by-one-digit practice answers in the\n# table with the newly generated values, matching each unique original\n# string and swapping it for its replacement.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old=\"64\u00f75=12, 4\"; New=\"67\u00f79=7, 4\"},\n    @{Old=\"57\u00f72=28, 1\"; New=\"45\u00f72=22, 1\"},\n    @{Old=\"65\u00f79=7, 2\"; New=\"88\u00f79=9, 7\"},\n    @{Old=\"62\u00f79=6, 8\"; New=\"50\u00f78=6, 2\"},\n    @{Old=\"92\u00f74=23, 0\"; New=\"23\u00f76=3, 5\"},\n    @{Old=\"60\u00f74=15, 0\"; New=\"51\u00f77=7, 2\"},\n    @{Old=\"78\u00f77=11, 1\"; New=\"83\u00f76=13, 5\"},\n    @{Old=\"79\u00f77=11, 2\"; New=\"38\u00f72=19, 0\"},\n    @{Old=\"49\u00f73=16, 1\"; New=\"19\u00f72=9, 1\"},\n    @{Old=\"78\u00f79=8, 6\"; New=\"28\u00f73=9, 1\"},\n    @{Old=\"93\u00f78=11, 5\"; New=\"81\u00f77=11, 4\"},\n    @{Old=\"76\u00f75=15, 1\"; New=\"49\u00f77=7, 0\"},\n    @{Old=\"45\u00f75=9, 0\"; New=\"29\u00f75=5, 4\"},\n    @{Old=\"36\u00f78=4, 4\"; New=\"97\u00f79=10, 7\"},\n    @{Old=\"25\u00f78=3, 1\"; New=\"76\u00f72=38, 0\"},\n    @{Old=\"66\u00f74=16, 2\"; New=\"46\u00f76=7, 4\"},\n    @{Old=\"42\u00f74=10, 2\"; New=\"32\u00f73=10, 2\"},\n    @{Old=\"11\u00f75=2, 1\"; New=\"45\u00f78=5, 5\"},\n    @{Old=\"67\u00f75=13, 2\"; New=\"55\u00f78=6, 7\"},\n    @{Old=\"48\u00f78=6, 0\"; New=\"78\u00f78=9, 6\"},\n    @{Old=\"43\u00f76=7, 1\"; New=\"40\u00f79=4, 4\"},\n    @{Old=\"68\u00f72=34, 0\"; New=\"66\u00f75=13, 1\"},\n    @{Old=\"69\u00f72=34, 1\"; New=\"52\u00f77=7, 3\"},\n    @{Old=\"62\u00f72=31, 0\"; New=\"61\u00f79=6, 7\"},\n    @{Old=\"29\u00f78=3, 5\"; New=\"18\u00f75=3, 3\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.New\n    $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
